$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-detected by Excel as numbers, so they remain text like the originals.
$textFormatCells = @("D5", "D6", "D9", "D11", "D14", "D16", "D19", "D20", "D21", "D22", "D25", "D32", "D35", "D37", "D41", "D43", "D47", "D50", "D51")
foreach ($c in $textFormatCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated cell values from the crypto price refresh
$ws.Range("D2").Value = "64.871.72"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "3.554.32"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "597.68"
$ws.Range("D6").Value = "134.95"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("D7").Value = "3.553.28"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").Value = "6.93"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "4.157.70"
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.560.42"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "26.83"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "64.986.10"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "9.91"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").Value = "14.29"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("D21").Value = "5.78"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "387.53"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +3.58%  "
$ws.Range("D24").Value = "3.701.37"
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("D25").Value = "73.77"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("E28").Value = "  +4.04%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("D32").Value = "1.48"
$ws.Range("E32").Value = "  +24.16%  "
$ws.Range("D33").Value = "3.552.54"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D35").Value = "23.90"
$ws.Range("E35").Value = "  +3.83%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "168.24"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("D41").Value = "0.0800"
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").Value = "26.76"
$ws.Range("E43").Value = "  +7.40%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").Value = "1.19"
$ws.Range("E47").Value = "  +3.39%  "
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").Value = "2.471.93"
$ws.Range("E49").Value = "  +11.70%  "
$ws.Range("D50").Value = "6.87"
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").Value = "0.863"
$ws.Range("E51").Value = "  +8.12%  "
